$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "notes" column header in D2 (matches header style of C2) ---
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = "notes"

# --- Add "First note" to D3, D4, D5 (matching each row's column-C style) ---
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "First note"

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "First note"

$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = "First note"

# --- Add new row 6: duplicate of row 3 data (shashwat / 123 / Kalkaji) plus "Second note" ---
# Copy formatting for A6:D6 from the row above (row 5), which carries the repeating data-row style.
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)

$ws.Range("A6").Value = "shashwat"
$ws.Range("B6").Value = 123
$ws.Range("C6").Value = "Kalkaji"
$ws.Range("D6").Value = "Second note"
